$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 2,3,5,6 (GRUPO EGIDO PINTOBASKET vs LUJISA GUADALAJARA BASKET) ---
$ws.Range("D2").Value = "S. RODRIGUEZ GREGORIO"
$ws.Range("D3").Value = "J. AYALA AYLLÓN"
$ws.Range("D5").Value = "M. DEL VALLE REGALADO"
$ws.Range("D6").Value = "M. BATLLE BERNARDO"

# --- Rows 9,10,11 (GRUPO EGIDO PINTOBASKET vs LUJISA GUADALAJARA BASKET) ---
$ws.Range("D9").Value = "A. ARMSTRONG"
$ws.Range("D10").Value = "L. VALERA VILLEGAS"
$ws.Range("D11").Value = "B. OJEDA OCHOA"

# --- Row 13 (C. D. MENSAJERO ISLA DE LA PALMA vs REAL CANOE N.C.) ---
$ws.Range("D13").Value = "P. RODRIGUEZ RIVERO"
$ws.Range("Q13").Value = 1
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 0
$ws.Range("T13").Value = 2
$ws.Range("U13").Value = 1
$ws.Range("V13").Value = 1
$ws.Range("W13").Value = 0

# --- Row 14 ---
$ws.Range("D14").Value = "A. APARICIO IZQUIERDO"
$ws.Range("M14").Value = 0
$ws.Range("N14").Value = 0
$ws.Range("P14").Value = ""
$ws.Range("R14").Value = 1
$ws.Range("S14").Value = 1
$ws.Range("T14").Value = 0
$ws.Range("V14").Value = 0
$ws.Range("W14").Value = 15.06

# --- Row 15 ---
$ws.Range("D15").Value = "I. REBERGEN"
$ws.Range("M15").Value = 2
$ws.Range("N15").Value = 2
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 0
$ws.Range("W15").Value = 13.25

# --- Row 34 (BALONCESTO TELDE vs ADC BOADILLA) ---
$ws.Range("D34").Value = "J. JIMENEZ HERNANDEZ"
$ws.Range("Q34").Value = 1
$ws.Range("T34").Value = 1
$ws.Range("U34").Value = 1

# --- Row 35 ---
$ws.Range("D35").Value = "A. SANTANA OJEDA"
$ws.Range("S35").Value = 0

# --- Row 36 ---
$ws.Range("D36").Value = "A. MARTIN GONZALEZ"
$ws.Range("Q36").Value = 0
$ws.Range("S36").Value = 1
$ws.Range("T36").Value = 0
$ws.Range("U36").Value = 0
